# Apply text replacements to the three-digit division problems table.
$d = $word.ActiveDocument

$replacements = @(
    @("296÷7=", "121÷4="),
    @("261÷4=", "841÷3="),
    @("331÷7=", "376÷9="),
    @("946÷6=", "101÷4="),
    @("438÷8=", "198÷8="),
    @("684÷3=", "111÷8="),
    @("137÷4=", "377÷3="),
    @("911÷8=", "756÷2="),
    @("346÷5=", "101÷8="),
    @("919÷2=", "823÷7="),
    @("954÷9=", "317÷7="),
    @("970÷4=", "350÷3="),
    @("135÷2=", "519÷5="),
    @("719÷9=", "432÷7="),
    @("320÷8=", "571÷2="),
    @("477÷7=", "225÷6="),
    @("962÷6=", "106÷3="),
    @("777÷8=", "745÷8="),
    @("956÷2=", "450÷3="),
    @("625÷6=", "737÷5="),
    @("769÷3=", "557÷7="),
    @("486÷4=", "439÷7="),
    @("752÷7=", "460÷7="),
    @("507÷9=", "606÷3="),
    @("176÷9=", "522÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
